$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestSubject")
$tbl = $ws.ListObjects.Item(1)

# 1. Add 5 new table columns at the end (table grows from 9 to 14 columns)
for ($i = 0; $i -lt 5; $i++) {
  $tbl.ListColumns.Add() | Out-Null
}

# 2. Copy the "last column" border/format from I (old last column, "Position during
#    measurement") into N (new last column) before we repurpose I for other content
$ws.Range("I6:I7").Copy()
$ws.Range("N6:N7").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# 3. Copy the "middle column" border/format from H into I, J, K, L, M so every new
#    header/data cell gets the same thin/medium table borders as its neighbours
$ws.Range("H6:H7").Copy()
$ws.Range("I6:M7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# 4. M6 (now second-to-last header) loses its right border since N6 supplies the
#    separating border on its own left edge
$ws.Range("M6").Borders.Item(10).LineStyle = -4142   # xlEdgeRight -> none

# 5. Rename headers to their final text, left-to-right, mirroring how a person typing
#    across the row would create the new shared strings in that order
$tbl.ListColumns.Item(7).Range.Cells.Item(1,1).Value  = "Stimulation current [mA]"
$tbl.ListColumns.Item(8).Range.Cells.Item(1,1).Value  = "Stimulation frequency [Hz]"
$tbl.ListColumns.Item(9).Range.Cells.Item(1,1).Value  = "Stimulation pulse width [us]"
$tbl.ListColumns.Item(10).Range.Cells.Item(1,1).Value = "Stimulation time [s]"
$tbl.ListColumns.Item(11).Range.Cells.Item(1,1).Value = "Rest time (fatigue) [s]"
$tbl.ListColumns.Item(12).Range.Cells.Item(1,1).Value = "Fatigue repetitions"
$tbl.ListColumns.Item(13).Range.Cells.Item(1,1).Value = "Leg side"
$tbl.ListColumns.Item(14).Range.Cells.Item(1,1).Value = "Position during measurement"

# 6. Give the new columns sensible widths (narrow, uniform for the new stim-parameter
#    / "Leg side" columns, a bit wider for the final "Position during measurement")
$ws.Range("G1:M1").ColumnWidth = 13.29
$ws.Range("N1").ColumnWidth = 18.71

# 7. Scroll the view roughly to where the new columns were added and leave the
#    selection on the newly typed header, matching the saved view state
$win = $excel.ActiveWindow
$win.ScrollColumn = 5
$ws.Range("M7").Select()
